$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the stray column-level style on column A (it was never really used by
# any cell) while preserving the header cell's own formatting and the
# existing black-font style used by the data rows below it.
$headerValue = $ws.Range("A1").Value
$ws.Columns("A").ClearFormats()
$ws.Range("B1").Copy($ws.Range("A1"))
$ws.Range("A1").Value = $headerValue
$ws.Range("A2:A12").Font.Color = 0

# Add the new row of data (row 13)
$ws.Range("A13").Value = 21682000
$ws.Range("A13").Font.Color = 0
$ws.Range("B13").Value = "Kho Giao Hàng Nặng Hà Tĩnh"
$ws.Range("F13").Value = "Ca chiều"

# Match the selection shown in the saved workbook
$ws.Range("B22").Select()

# Page setup: A4 paper, portrait orientation
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
